{"js": "// Replace the date line and each \"A\u00d7B=\" problem text in the practice\n// sheet with the updated values from the new day's worksheet.\n// Each old string is unique within the document body, so a targeted\n// search + replace is used for every pair (search results are re-queried\n// each time so earlier replacements can't shift later matches).\n\nconst replacements = [\n  [\"2025-06-13 Friday\", \"2025-06-14 Saturday\"],\n  [\"564\u00d79=\", \"536\u00d74=\"],\n  [\"106\u00d73=\", \"698\u00d75=\"],\n  [\"711\u00d79=\", \"331\u00d74=\"],\n  [\"857\u00d72=\", \"842\u00d78=\"],\n  [\"850\u00d74=\", \"902\u00d79=\"],\n  [\"113\u00d75=\", \"357\u00d76=\"],\n  [\"993\u00d77=\", \"210\u00d73=\"],\n  [\"367\u00d79=\", \"870\u00d77=\"],\n  [\"212\u00d77=\", \"213\u00d77=\"],\n  [\"970\u00d77=\", \"792\u00d75=\"],\n  [\"657\u00d79=\", \"358\u00d73=\"],\n  [\"532\u00d76=\", \"989\u00d79=\"],\n  [\"779\u00d74=\", \"944\u00d78=\"],\n  [\"239\u00d72=\", \"507\u00d73=\"],\n  [\"545\u00d73=\", \"920\u00d78=\"],\n  [\"906\u00d73=\", \"652\u00d78=\"],\n  [\"162\u00d76=\", \"179\u00d72=\"],\n  [\"504\u00d72=\", \"262\u00d74=\"],\n  [\"579\u00d73=\", \"477\u00d75=\"],\n  [\"578\u00d77=\", \"736\u00d72=\"],\n  [\"444\u00d74=\", \"131\u00d77=\"],\n  [\"660\u00d76=\", \"426\u00d73=\"],\n  [\"240\u00d79=\", \"314\u00d77=\"],\n  [\"440\u00d77=\", \"192\u00d77=\"],\n  [\"899\u00d78=\", \"652\u00d78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"A\u00d7B=\" problem text in the practice\n# sheet with the updated values from the new day's worksheet.\n# Each old string is unique within the document, so Find/Replace against\n# the whole document Range is safe for every pair.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-13 Friday\", \"2025-06-14 Saturday\"),\n    @(\"564\u00d79=\", \"536\u00d74=\"),\n    @(\"106\u00d73=\", \"698\u00d75=\"),\n    @(\"711\u00d79=\", \"331\u00d74=\"),\n    @(\"857\u00d72=\", \"842\u00d78=\"),\n    @(\"850\u00d74=\", \"902\u00d79=\"),\n    @(\"113\u00d75=\", \"357\u00d76=\"),\n    @(\"993\u00d77=\", \"210\u00d73=\"),\n    @(\"367\u00d79=\", \"870\u00d77=\"),\n    @(\"212\u00d77=\", \"213\u00d77=\"),\n    @(\"970\u00d77=\", \"792\u00d75=\"),\n    @(\"657\u00d79=\", \"358\u00d73=\"),\n    @(\"532\u00d76=\", \"989\u00d79=\"),\n    @(\"779\u00d74=\", \"944\u00d78=\"),\n    @(\"239\u00d72=\", \"507\u00d73=\"),\n    @(\"545\u00d73=\", \"920\u00d78=\"),\n    @(\"906\u00d73=\", \"652\u00d78=\"),\n    @(\"162\u00d76=\", \"179\u00d72=\"),\n    @(\"504\u00d72=\", \"262\u00d74=\"),\n    @(\"579\u00d73=\", \"477\u00d75=\"),\n    @(\"578\u00d77=\", \"736\u00d72=\"),\n    @(\"444\u00d74=\", \"131\u00d77=\"),\n    @(\"660\u00d76=\", \"426\u00d73=\"),\n    @(\"240\u00d79=\", \"314\u00d77=\"),\n    @(\"440\u00d77=\", \"192\u00d77=\"),\n    @(\"899\u00d78=\", \"652\u00d78=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
